# Add two new bold bullet points right after the
# "I didn't actually write 0x02 to POWER_CTL :/" bullet:
#   - "Not sure why temperature sensor is always 25 now."
#   - "Never mind, it's actually 25."

$d = $word.ActiveDocument

$anchorText = "I didn’t actually write 0x02 to POWER_CTL :/"

# Locate the anchor paragraph by exact text (ignoring the trailing
# paragraph-mark character Word appends to Paragraph.Range.Text).
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq 0) {
    Write-Output "Anchor paragraph not found!"
} else {
    $anchorPara = $d.Paragraphs.Item($anchorIndex)

    # Inserting a paragraph after the anchor clones its pPr (ListBullet
    # style, numId=3 bullet list) onto the new, still-empty paragraph.
    $anchorPara.Range.InsertParagraphAfter()

    $p1 = $d.Paragraphs.Item($anchorIndex + 1)
    $p1.Range.Text = "Not sure why temperature sensor is always 25 now."
    $r1 = $p1.Range.Duplicate
    $r1.MoveEnd(1, -1)   # exclude the paragraph mark so <w:pPr> stays clean
    $r1.Bold = 1

    $p1.Range.InsertParagraphAfter()

    $p2 = $d.Paragraphs.Item($anchorIndex + 2)
    $p2.Range.Text = "Never mind, it’s actually 25."
    $r2 = $p2.Range.Duplicate
    $r2.MoveEnd(1, -1)
    $r2.Bold = 1

    Write-Output "Inserted 2 paragraphs after paragraph $anchorIndex."
}
